# Data.xlsx re-upload: refreshed column B ("total") figures for rows 2-38
# (except row 31, which is unchanged) and cleared the direct cell format
# that had been applied to column C ("Delta") on those same rows, plus
# moved the saved cell selection to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column B values keyed by row number.
$newTotals = @{
    2  = 25515
    3  = 13935
    4  = 39090
    5  = 35745
    6  = 30600
    7  = 28211
    8  = 26646
    9  = 23332
    10 = 12010
    11 = 16306
    12 = 34647
    13 = 28550
    14 = 22225
    15 = 19359
    16 = 17584
    17 = 8921
    18 = 15365
    19 = 30439
    20 = 26858
    21 = 23511
    22 = 22771
    23 = 21366
    24 = 10940
    25 = 33128
    26 = 29413
    27 = 25312
    28 = 24668
    29 = 26842
    30 = 26167
    32 = 43662
    33 = 40731
    34 = 34914
    35 = 35852
    36 = 37251
    37 = 34444
    38 = 18444
}

foreach ($row in 2..38) {
    if ($newTotals.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value = $newTotals[$row]
    }
    # Clear the direct formatting that had been applied to column C
    # ("Delta"), reverting those cells to the workbook's default style.
    $ws.Cells.Item($row, 3).ClearFormats()
}

# Restore the saved selection to F7, as recorded in the sheet view.
[void]$ws.Range("F7").Select()
